$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
        @(1101.58064516129,1161.321428571429,1155.41935483871,1112.866666666667,883.3387096774194,1014.316666666667,1016.774193548387,980.0967741935484,1010.433333333333,983.6612903225806,1086.516666666667,1112.564516129032),
        @(1104.838709677419,1157.803571428571,1154.41935483871,1110.183333333333,875.5645161290323,1014.016666666667,1016.58064516129,981.4193548387096,1010.666666666667,981.6129032258065,1084.208333333333,1108.516129032258),
        @(1104.629032258065,1154.142857142857,1151.85,1108.3,869.9354838709677,1009.233333333333,1014.564516129032,981.6935483870968,1008.766666666667,981.6935483870968,1087.85,1109.177419354839),
        @(1102.209677419355,1155.875,1154.629032258065,1109.566666666667,870.2096774193549,1009.183333333333,1015.064516129032,981.516129032258,1010.116666666667,983.4677419354839,1088.516666666667,1108.661290322581),
        @(1102.145161290323,1153.464285714286,1146.403225806452,1106.116666666667,869.8548387096774,1007.066666666667,1014.951612903226,980.8548387096774,1009.85,983.5,1086.583333333333,1108.854838709677),
        @(1102.822580645161,1155.821428571429,1149.064516129032,1105.516666666667,878.741935483871,1005.166666666667,1014.790322580645,981.483870967742,1010.483333333333,979.5483870967741,1090.15,1110.403225806452),
        @(1115.306451612903,1163,1162.403225806452,1112.633333333333,888.0806451612904,1013.85,1017.838709677419,986.2258064516129,1017.05,996.4193548387096,1101.116666666667,1126.241935483871),
        @(1152.209677419355,1182.535714285714,1184.290322580645,1128.85,912.4516129032259,1041.833333333333,1032.838709677419,1007.58064516129,1031.383333333333,1012.887096774194,1118.716666666667,1141.032258064516),
        @(1197.403225806452,1223.678571428571,1238.193548387097,1181.6,980.4677419354839,1118.266666666667,1097.548387096774,1050.209677419355,1071.433333333333,1055.854838709677,1173.533333333333,1185.645161290323),
        @(1256.677419354839,1307.660714285714,1288.725806451613,1245.666666666667,1037.370967741935,1164.316666666667,1145.322580645161,1080.016129032258,1133.866666666667,1120.483870967742,1212.166666666667,1238.193548387097),
        @(1306.677419354839,1319.107142857143,1316.596774193548,1284.55,1072.096774193548,1195.166666666667,1163.967741935484,1115.290322580645,1177.833333333333,1185.983870967742,1240.25,1255.741935483871),
        @(1289.225806451613,1335.964285714286,1324.806451612903,1303.083333333333,1081,1212.266666666667,1182.935483870968,1123.967741935484,1197.716666666667,1203.548387096774,1252.466666666667,1259.806451612903),
        @(1290.645161290323,1337.696428571429,1334.661290322581,1325.716666666667,1100.790322580645,1211.066666666667,1184.645161290323,1127.322580645161,1202.85,1204.629032258065,1259.066666666667,1273.161290322581),
        @(1300.048387096774,1356.214285714286,1344.177419354839,1310,1103.451612903226,1231.4,1186.612903225806,1131.58064516129,1203.95,1195.338709677419,1255.033333333333,1267.91935483871),
        @(1278.41935483871,1331.035714285714,1317.435483870968,1297.4,1076.822580645161,1212.033333333333,1174,1121.354838709677,1186.066666666667,1180.193548387097,1246.05,1252.564516129032),
        @(1262.838709677419,1301.607142857143,1295.516129032258,1272.933333333333,1049.725806451613,1167.25,1138.677419354839,1100.435483870968,1174.716666666667,1161.112903225806,1234.7,1235.322580645161),
        @(1229.145161290323,1265.25,1264.612903225806,1225.833333333333,998.483870967742,1119.016666666667,1096.854838709677,1061.5,1125.583333333333,1109.532258064516,1203.433333333333,1211.564516129032),
        @(1203.258064516129,1244.339285714286,1243.112903225806,1197.733333333333,966.1612903225806,1100.366666666667,1079.91935483871,1036.903225806452,1094.95,1075.129032258065,1182.65,1191.483870967742),
        @(1184.58064516129,1224.375,1223.870967741935,1174.833333333333,946.6290322580645,1077.633333333333,1060.645161290323,1019.274193548387,1072.55,1050.274193548387,1159.433333333333,1176.41935483871),
        @(1167.322580645161,1204.321428571429,1212.983870967742,1162.116666666667,936.741935483871,1060.116666666667,1044.967741935484,1008.016129032258,1054.8,1030.532258064516,1139.216666666667,1160.758064516129),
        @(1153.354838709677,1193.125,1195.709677419355,1150.133333333333,923.241935483871,1050.416666666667,1031.903225806452,1000.451612903226,1041.733333333333,1014.209677419355,1127.683333333333,1153.322580645161),
        @(1138.451612903226,1178.625,1166.338709677419,1131.466666666667,905.7096774193549,1036.283333333333,1022.08064516129,992.516129032258,1031,1001.096774193548,1108.133333333333,1143.096774193548),
        @(1115.661290322581,1172.107142857143,1159.064516129032,1116.6,893.3225806451613,1024.533333333333,1021.725806451613,988.6612903225806,1025.616666666667,998.4677419354839,1106.85,1135.370967741935),
        @(1106.935483870968,1165.017857142857,1157.838709677419,1112.383333333333,884.258064516129,1014.416666666667,1016.758064516129,985.0645161290323,1016.616666666667,991.6290322580645,1094.6,1117.306451612903)
)

for ($r = 0; $r -lt 24; $r++) {
    for ($c = 0; $c -lt 12; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $data[$r][$c]
    }
}
